$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.237.61'
$ws.Range("E2").Value = '  -9.98%  '
$ws.Range("D3").Value = '2.386.37'
$ws.Range("E3").Value = '  -12.63%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '458.31'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -9.74%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '129.24'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -8.67%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.480'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -10.19%  '
$ws.Range("D9").Value = '2.398.90'
$ws.Range("E9").Value = '  -12.67%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0941'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -10.05%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.24'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -14.43%  '
$ws.Range("E12").Value = '  -11.11%  '
$ws.Range("E13").Value = '  -4.61%  '
$ws.Range("D14").Value = '2.794.93'
$ws.Range("E14").Value = '  -12.88%  '
$ws.Range("D15").Value = '53.265.89'
$ws.Range("E15").Value = '  -9.79%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '19.46'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -10.94%  '
$ws.Range("E17").Value = '  -5.20%  '
$ws.Range("D18").Value = '2.383.20'
$ws.Range("E18").Value = '  -12.70%  '
$ws.Range("E19").Value = '  -13.25%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '304.88'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -11.91%  '
$ws.Range("E21").Value = '  -15.77%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("E24").Value = '  -15.64%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '55.49'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -12.30%  '
$ws.Range("E26").Value = '  -1.23%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.381'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -10.48%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.151'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -13.03%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.465.42'
$ws.Range("E29").Value = '  -13.26%  '
$ws.Range("E30").Value = '  -6.75%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").Value = '0.0₃0716'
$ws.Range("E32").Value = '  -14.70%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '144.74'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.91%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '17.56'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -8.30%  '
$ws.Range("E35").Value = '  -13.23%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.95'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -8.37%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.49'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -17.37%  '
$ws.Range("E38").Value = '  -8.11%  '
$ws.Range("E39").Value = '  -16.93%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.993'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.22%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '32.77'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -9.51%  '
$ws.Range("E42").Value = '  -3.01%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.24'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -8.50%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0520'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -7.00%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '10.10'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("E46").Value = '  -11.99%  '
$ws.Range("D47").Value = '1.932.76'
$ws.Range("E47").Value = '  -11.66%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0216'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -5.93%  '
$ws.Range("E49").Value = '  -3.22%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '4.18'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -13.39%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '16.30'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -14.85%  '
